$d = $word.ActiveDocument

# Minimal WordOpenXML package fragment containing a single run with the
# "(Flashback?) " text, using the sz/szCs=28 run formatting used
# throughout this document's bullet list. Using Range.InsertXML (rather
# than setting Range.Text) makes the inserted text become its own
# <w:r> element instead of being merged into the neighboring run.
$flashbackRunXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">(Flashback?) </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

function Insert-FlashbackTagBefore($searchText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)
    if ($rng.Find.Found) {
        $insertAt = $d.Range($rng.Start, $rng.Start)
        $insertAt.InsertXML($flashbackRunXml)
    }
}

# 1) "Zo's favorite sibling and best friend is conscripted as a miner..."
#    Split the existing run so the paragraph starts with a new
#    "(Flashback?) " run followed by the original, unmodified run.
Insert-FlashbackTagBefore("Zo" + [char]0x2019 + "s favorite sibling and best friend is conscripted")

# 2) "Zo announces to their parents that they want to study..."
#    Insert a brand-new "(Flashback?) " run as the first run of the
#    paragraph, immediately followed by the existing run.
Insert-FlashbackTagBefore("Zo announces to their parents")
